# Apply the cryptocurrency price/volume refresh described in the commit.
# Column D ("Price") values are forced to remain text (matching the original
# inlineStr cells) by temporarily setting a Text number format before writing
# the value, then restoring the default "Normal" style so no stray style index
# is left attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '23.271.05'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -2.63%  '

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.600.46'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -3.36%  '

# Row 4
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.004'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  +0.17%  '

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '1.004'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +0.21%  '

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '302.53'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -2.17%  '

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.3777'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -2.94%  '

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.3683'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -4.14%  '

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '49.48'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -2.88%  '

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '1.003'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +0.13%  '

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '1.282'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -5.39%  '

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.08130'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -4.01%  '

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '23.01'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -3.84%  '

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '6.670'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -6.80%  '

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '7.597'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -3.70%  '

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '0.00001271'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -2.66%  '

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '1.597.40'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -3.35%  '

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '91.61'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -3.37%  '

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.06831'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -2.63%  '

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '18.59'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -6.19%  '

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '6.628'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -4.06%  '

# Row 22
$ws.Range("E22").Value = '  +0.24%  '

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '13.21'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -2.82%  '

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '23.277.71'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -2.59%  '

# Row 25
$ws.Range("B25").Value = 'LidoDAOToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.987'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -1.91%  '

# Row 26
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '2.365'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -5.24%  '

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '21.20'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -3.94%  '

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '151.10'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -1.11%  '

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '5.311'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -1.78%  '

# Row 30
$ws.Range("E30").Value = '  -4.88%  '

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '2.474'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -0.83%  '

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '7.166'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -7.88%  '

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '1.774.61'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -3.28%  '

# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.9708'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -6.12%  '

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.07753'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -3.80%  '

# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.02794'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -5.74%  '

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '6.338'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -4.94%  '

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '10.30'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -6.37%  '

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.2561'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -4.59%  '

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.08880'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -2.61%  '

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '1.393'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -1.90%  '

# Row 42
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.7202'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -4.48%  '

# Row 43
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '12.88'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -4.56%  '

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '16.24'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -0.07%  '

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.6649'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -4.34%  '

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.326'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -5.59%  '

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '1.003'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +0.16%  '

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '3.983'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -2.30%  '

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.08009'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -2.98%  '

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '132.08'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -1.79%  '

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '1.182'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -4.24%  '
